$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value to 400 (was 200)
$ws.Range("A2").Value = "param_pv1_area"
$ws.Range("B2").Value = 400

# Delete rows 3 through 10 (old params) entirely
$ws.Range("A3:B10").EntireRow.Delete()
